# Apply updated scoring data: adds an option to configure the cost function
# so that it can vary for prelims / finals / prelim-finals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scoringData = @{
  2 = @(0,0,4,16,0,0,10,13,0,1,0,7,0,0,0,0,3.25,0)
  3 = @(12,4,7,11,0,0,10,5,11,9,11,7,3,1,3,0,1.25,2.75)
  4 = @(12,10,1,10,16,16,10,5,7,9,5,12,3,2.5,3,4,1.25,1.75)
  5 = @(0,0,0,11,0,0,0,0,2,0,0,0,0,0,0,0,0,0.5)
  6 = @(0,0,0,10,0,2,0,0,1,0,5,0,0,0,0,0,0,0.25)
  7 = @(0,0,4,16,0,0,10,5,0,0,3,12,0,0,0,0,1.25,0)
  8 = @(1,7,13,16,3,0,0,0,10,13,7,12,0.25,1.75,0.25,0.75,0,2.5)
  9 = @(13,2,0,2,11,0,0,0,0,0,0,0,3.25,0.5,3.25,2.75,0,0)
  10 = @(0,0,0,7,9,1,0,10,7,0,0,0,0,0,0,2.25,2.5,1.75)
  11 = @(13,16,16,16,16,16,16,16,16,16,16,16,3.25,4,3.25,4,4,4)
  12 = @(0,0,9,12,0,0,0,0,0,0,0,2,0,0,0,0,0,0)
  13 = @(7,0,0,7,0,0,0,0,0,9,4,0,1.75,0,1.75,0,0,0)
}

foreach ($rowNum in $scoringData.Keys) {
  $rowValues = $scoringData[$rowNum]
  $rowArray = New-Object "object[,]" 1,18
  for ($col = 0; $col -lt 18; $col++) { $rowArray[0,$col] = $rowValues[$col] }
  $targetRange = $ws.Range("B" + $rowNum + ":S" + $rowNum)
  $targetRange.Value = $rowArray
}

# Widen the relay columns (N:Q) which now hold longer computed values
$ws.Columns.Item(14).ColumnWidth = 10.666666666666666
$ws.Columns.Item(15).ColumnWidth = 11.666666666666666
$ws.Columns.Item(16).ColumnWidth = 13.998697916666666
$ws.Columns.Item(17).ColumnWidth = 14.166666666666666

# Restore the view state (zoom level and selection) used while reviewing the data
$excel.ActiveWindow.Zoom = 88
$ws.Range("A1:S13").Select() | Out-Null
